$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet/tab (device_classification_fields -> nsde_fields)
$ws.Name = "nsde_fields"

# Row 12: inactivation_date
$ws.Range("B12").Value = "inactivation_date"
$ws.Range("B12").ClearFormats()
$ws.Range("C12").Value = "string"
$ws.Range("D12").Value = "The date on which registration or listing data was inactivated by FDA due to inaccuracies, incompleteness or incompliance."

# Row 13: reactivation_date
$ws.Range("B13").Value = "reactivation_date"
$ws.Range("C13").Value = "string"
$ws.Range("D13").Value = "The date on which a previously FDA inactivated registration or listing data is reactivated."

# Row heights to match wrapped-text content
$ws.Rows.Item(12).RowHeight = 34
$ws.Rows.Item(13).RowHeight = 17

# Zoom + selection to match author's final view state
$excel.ActiveWindow.Zoom = 140
$ws.Range("B12").Select() | Out-Null
